# Mediators_Upload_Sample_Document.xlsx - sample doc update for interpreter bulk updates
#
#  - header row: "firstName"/"lastName" -> "first_name"/"last_name"
#  - P2:Q2 (the two boolean TRUE() columns) pick up the same cell format as the
#    rest of the header/data range instead of their own separate (duplicate)
#    style record
#  - active selection moves to B1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the first two header columns.
$ws.Range("A1").Value = "first_name"
$ws.Range("B1").Value = "last_name"

# P2:Q2 were carrying their own (duplicate) cell format; re-apply the same
# number format already in effect so the engine folds them back onto the
# shared style used by the rest of the row instead of a standalone one.
$ws.Range("P2:Q2").NumberFormat = $ws.Range("A1").NumberFormat

# Move the active selection to B1.
[void]$ws.Range("B1").Select()
